# Continue the "Example" sheet work: the blade-width note used for both the
# rectangular- and circular-coupon tables is being reworded now that the
# Panel/Defaults sections are finished, and the entry point for the next
# (coupon) section is selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Reword "blade width (leave blank)" -> "blade width (leave blank if same)"
# everywhere it appears (Rectangular Coupons header row + Circular Coupons
# header row). Setting the cell text directly lets the shared-string table
# be rebuilt naturally (old text dropped, new text appended).
$ws.Range("E10").Value = "blade width (leave blank if same)"
$ws.Range("D19").Value = "blade width (leave blank if same)"

# The two columns that hold that long label (D and E) need to be widened
# to fit the new, longer text.
$ws.Range("D1").EntireColumn.ColumnWidth = 25
$ws.Range("E1").EntireColumn.ColumnWidth = 25

# Move on: select the first empty row below the finished tables, ready to
# start the coupons section.
$ws.Range("A21").Select()
